$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Balance" header column
$ws.Range("F1").Value = "Balance"

# Row 2 - bicycle purchase (expense)
$ws.Range("A2").Value = "15:36"
$ws.Range("B2").Value = "buying a bicycle"
$ws.Range("C2").Value = "purchase"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'200"
$ws.Range("F2").Value = -200

# Row 3 - paint sale (income)
$ws.Range("A3").Value = "15:37"
$ws.Range("B3").Value = "selling paint"
$ws.Range("C3").Value = "sell"
$ws.Range("D3").Value = "'200"
$ws.Range("E3").Value = "'0"
$ws.Range("F3").Value = 0
